# Generate Report for Handback
# Refresh the handback-status report: update the "Latest HO Xliff Generate
# Date" / handoff / handback timestamps for the 25cb016f-... file now that a
# new handback round has completed.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 25cb016f-... row (row 2) advances to the new generation timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 08:52:16"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback
# DateTime (K) for the 25cb016f-... row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 08:52:12"
$wsZhCn.Range("K2").Value = "2016-09-03 08:52:29"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback
# DateTime (K) for the 25cb016f-... row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 08:52:16"
$wsDeDe.Range("K2").Value = "2016-09-03 08:52:36"
